$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on Price column (D) cells that will receive numeric-looking
# strings, so Excel does not auto-convert them to real numbers.
$priceCells = @(
    "D2", "D3", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16",
    "D17", "D18", "D19", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31",
    "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D42", "D43", "D44", "D45", "D46",
    "D50", "D51"
)
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '29.481.51'

# Row 3
$ws.Range("D3").Value = '1.879.78'
$ws.Range("E3").Value = '  +1.15%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").Value = '0.7159'
$ws.Range("E5").Value = '  +1.20%  '

# Row 6
$ws.Range("D6").Value = '241.92'
$ws.Range("E6").Value = '  +1.55%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").Value = '0.07864'
$ws.Range("E8").Value = '  -2.01%  '

# Row 9
$ws.Range("D9").Value = '0.3124'
$ws.Range("E9").Value = '  +3.02%  '

# Row 10
$ws.Range("D10").Value = '25.28'
$ws.Range("E10").Value = '  +7.33%  '

# Row 11
$ws.Range("D11").Value = '0.08262'
$ws.Range("E11").Value = '  +0.90%  '

# Row 12
$ws.Range("D12").Value = '0.7320'
$ws.Range("E12").Value = '  +3.52%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.868.73'
$ws.Range("E13").Value = '  +0.36%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.288'
$ws.Range("E14").Value = '  +1.65%  '

# Row 15
$ws.Range("D15").Value = '91.27'
$ws.Range("E15").Value = '  +1.70%  '

# Row 16
$ws.Range("D16").Value = '29.429.03'
$ws.Range("E16").Value = '  +0.59%  '

# Row 17
$ws.Range("D17").Value = '5.941'
$ws.Range("E17").Value = '  +2.00%  '

# Row 18
$ws.Range("D18").Value = '248.47'
$ws.Range("E18").Value = '  +4.16%  '

# Row 19
$ws.Range("D19").Value = '0.000007898'
$ws.Range("E19").Value = '  -0.46%  '

# Row 20
$ws.Range("E20").Value = '  +0.34%  '

# Row 21
$ws.Range("D21").Value = '0.9994'
$ws.Range("E21").Value = '  +0.09%  '

# Row 22
$ws.Range("D22").Value = '7.948'
$ws.Range("E22").Value = '  +6.15%  '

# Row 23
$ws.Range("E23").Value = '  +0.07%  '

# Row 24
$ws.Range("D24").Value = '0.1589'
$ws.Range("E24").Value = '  +10.48%  '

# Row 25
$ws.Range("D25").Value = '164.18'
$ws.Range("E25").Value = '  +0.80%  '

# Row 26
$ws.Range("D26").Value = '9.046'
$ws.Range("E26").Value = '  +1.67%  '

# Row 27
$ws.Range("D27").Value = '18.33'
$ws.Range("E27").Value = '  +1.14%  '

# Row 28
$ws.Range("D28").Value = '1.361'
$ws.Range("E28").Value = '  -4.61%  '

# Row 29
$ws.Range("D29").Value = '1.497'
$ws.Range("E29").Value = '  +1.23%  '

# Row 30
$ws.Range("D30").Value = '4.375'
$ws.Range("E30").Value = '  +0.04%  '

# Row 31
$ws.Range("D31").Value = '4.132'
$ws.Range("E31").Value = '  +2.62%  '

# Row 32
$ws.Range("D32").Value = '0.05316'
$ws.Range("E32").Value = '  +2.38%  '

# Row 33
$ws.Range("D33").Value = '1.938'
$ws.Range("E33").Value = '  +0.41%  '

# Row 34
$ws.Range("D34").Value = '1.202'
$ws.Range("E34").Value = '  +3.54%  '

# Row 35
$ws.Range("D35").Value = '0.7241'
$ws.Range("E35").Value = '  +1.19%  '

# Row 36
$ws.Range("D36").Value = '2.679'
$ws.Range("E36").Value = '  +0.66%  '

# Row 37
$ws.Range("D37").Value = '0.01871'
$ws.Range("E37").Value = '  +0.69%  '

# Row 38
$ws.Range("D38").Value = '1.268.64'
$ws.Range("E38").Value = '  +11.47%  '

# Row 39
$ws.Range("E39").Value = '  +0.13%  '

# Row 40
$ws.Range("D40").Value = '0.9112'
$ws.Range("E40").Value = '  -2.76%  '

# Row 41
$ws.Range("E41").Value = '  +5.23%  '

# Row 42
$ws.Range("D42").Value = '6.107'
$ws.Range("E42").Value = '  +2.00%  '

# Row 43
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  +0.05%  '

# Row 44
$ws.Range("D44").Value = '103.75'
$ws.Range("E44").Value = '  +0.89%  '

# Row 45
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '0.5327'
$ws.Range("E45").Value = '  +0.56%  '

# Row 46
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '2.027.13'
$ws.Range("E46").Value = '  +0.99%  '

# Row 47
$ws.Range("E47").Value = '  +0.34%  '

# Row 48
$ws.Range("E48").Value = '  +13.05%  '

# Row 49
$ws.Range("E49").Value = '  +0.07%  '

# Row 50
$ws.Range("D50").Value = '0.4335'
$ws.Range("E50").Value = '  +1.48%  '

# Row 51
$ws.Range("D51").Value = '9.284'
$ws.Range("E51").Value = '  +1.13%  '

# Reset style reference on Price cells back to default (Normal) now that the
# text values are stored, so no extra explicit style indices linger on them.
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
